$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the "Fecha" (date) column: the old values were a bogus
# far-future date serial; replace with the real request dates ---
$ws.Range("A2").Value = "5/4/2017"
$ws.Range("A3").Value = "6/12/2017"
$ws.Range("A4").Value = "3/16/2017"

# --- New columns: Id solicitante (K), Nombre solicitante (L),
# Estado de solicitud (M) -- tracking the "resolution" request state ---
$ws.Range("K1").Value = "Id solicitante"
$ws.Range("L1").Value = "Nombre solicitante"
$ws.Range("M1").Value = "Estado de solicitud"

# Row 2
$ws.Range("K2").Value = "'123"
$ws.Range("K2").Style = "Normal"
$ws.Range("L2").Value = "ana maria "
$ws.Range("M2").Value = "PROCESADA"

# Row 3
$ws.Range("K3").Value = "'123"
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").Value = "ana maria "
$ws.Range("M3").Value = "PROCESADA"

# Row 4
$ws.Range("K4").Value = "'123"
$ws.Range("K4").Style = "Normal"
$ws.Range("L4").Value = "ana maria "
$ws.Range("M4").Value = "CANCELADO"

# Row 5: a stray underline formatted (non-hyperlink colored) empty cell,
# left over from extending the table formatting down one more row
$ws.Range("M5").Font.Underline = 2

# View state: user ended up with M2 selected, scrolled right a bit
$ws.Range("M2").Select()
$excel.ActiveWindow.ScrollColumn = 3

# Print setup: orientation explicitly set to portrait
$ws.PageSetup.Orientation = 1
